$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row N gets the D/J/K/L/M/P values that originally belonged to row Map[N]
# This reflects a reshuffle of the weekly price records while keeping the
# rest of each row (Mercado, Region, etc.) fixed.

$origD = @{2=44175; 3=44537; 4=44210; 5=44200; 6=44893; 7=44638; 8=44907; 9=44895; 10=44883}
$origJ = @{2=1400; 3=800; 4=1450; 5=1500; 6=3300; 7=800; 8=2300; 9=200; 10=290}
$origK = @{2=1900; 3=1300; 4=1600; 5=1400; 6=1200; 7=2500; 8=900; 9=1200; 10=1400}
$origL = @{2=2000; 3=1400; 4=1700; 5=1500; 6=1300; 7=2800; 8=1000; 9=1300; 10=1500}
$origM = @{2=1950; 3=1350; 4=1650; 5=1450; 6=1261; 7=2650; 8=952; 9=1255; 10=1434}

$map = @{2=3; 3=4; 4=8; 5=7; 6=5; 7=10; 8=9; 9=2; 10=6}

foreach ($row in 2..10) {
    $src = $map[$row]
    $ws.Cells.Item($row, 4).Value = $origD[$src]   # D: Fecha
    $ws.Cells.Item($row, 10).Value = $origJ[$src]  # J: Volumen
    $ws.Cells.Item($row, 11).Value = $origK[$src]  # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $origL[$src]  # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $origM[$src]  # M: Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $origM[$src]  # P: Precio $/Kg (same as M)
}
